$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up existing row 104 ---
$ws.Range("A104").Value = 45471.2916666667
$ws.Range("E104").Value = 0.689999997615814

# --- Append new row 105 ---
$ws.Range("A105").Value = 45474.2916666667
$ws.Range("B105").Value = 140
$ws.Range("C105").Value = 0.689999997615814
$ws.Range("D105").Value = 0.689999997615814
$ws.Range("E105").Value = 0.689999997615814
$ws.Range("F105").Value = 0.689999997615814

# --- Append new row 106 ---
$ws.Range("A106").Value = 45475.4291782407
$ws.Range("B106").Value = 9400
$ws.Range("C106").Value = 0.704999983310699
$ws.Range("D106").Value = 0.680000007152557
$ws.Range("E106").Value = 0.689999997615814
$ws.Range("F106").Value = 0.699999988079071

# Column A on the new rows needs the same date/time display style as the
# rest of the date column (style index 1 -> numFmt "yyyy-mm-dd hh:mm:ss").
$ws.Range("A104").Copy()
$ws.Range("A105:A106").PasteSpecial(-4122)

# Column G ("adj_close") and H ("ticker") in this sheet are stored as *text*
# that happens to look like a number (a quirk of the R export), so we must
# force a Text number format before assigning, otherwise Excel will coerce
# the numeric-looking strings back into real numbers.
$ws.Range("G105:G106").NumberFormat = "@"
$ws.Range("G105").Value = "0.689999997615814"
$ws.Range("G106").Value = "0.699999988079071"
$ws.Range("H105").Value = "BWZ.MI"
$ws.Range("H106").Value = "BWZ.MI"

# Restore the default style on the text cells we just reformatted so they
# match the rest of the sheet (no explicit style index).
$ws.Range("G105:G106").Style = "Normal"
